# Release-Notes.xlsx update:
# A new "Folder Inventory" entry was logged before the prior top entry,
# shifting all existing rows down by one. Metadata/Summary counters and
# timestamps are refreshed accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Folder Inventory sheet: insert a new row 2 for the newest folder
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Folder Inventory")

$ws.Range("A2:E2").Insert(-4121)   # xlShiftDown
$ws.Range("A2:E2").ClearFormats()  # keep default (unstyled) look like other data rows

$ws.Range("A2").Value = "Developing_a_Custom_RAG_App_Using_Azure_AI_Foundry"
$ws.Range("B2").Value = "Developing_a_Custom_RAG_App_Using_Azure_AI_Foundry"
$ws.Range("C2").Value = "2025-06-13 13:08:47 +0530"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Root"

# ---------------------------------------------------------------------
# 2. Metadata sheet: refresh generated-on timestamp and counters
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2025-06-13 07:39:05 UTC"
$meta.Range("B4").Value = 75

# "Workflow Run" is stored as text ("4"), not a number - force text type
# without leaving a stray style behind.
$meta.Range("B5").NumberFormat = "@"
$meta.Range("B5").Value = "4"
$meta.Range("B5").ClearFormats()

# ---------------------------------------------------------------------
# 3. Summary sheet: refresh folder counts and most recent update time
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B2").Value = 75
$summary.Range("B3").Value = 75
$summary.Range("B5").Value = "2025-06-13 13:08:47 +0530"
